# "fix run without testdata file" — add the missing test-data row to the
# "Behandeling" sheet (row 3) that the run expected to find, and make that
# sheet the active one (tab + selection), matching the original author's
# commit.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Behandeling")

# New test-data row: Dag/Maand/Jaar_behandeling, Type, Behandelaar,
# Naam_Behandeling, Afdeling, Urgentie.
$ws1.Range("A3").Value = 19
$ws1.Range("B3").Value = 7
$ws1.Range("C3").Value = 2021
$ws1.Range("D3").Value = "Opname"
$ws1.Range("E3").Value = "Dokter anders"
$ws1.Range("F3").Value = "Knie operatie"
$ws1.Range("G3").Value = "Radiologie"
$ws1.Range("H3").Value = "Laag"

# Make "Behandeling" the active sheet/tab, with A4 selected (just below the
# newly-added row) instead of the previous "Patient" tab / B7 selection.
$ws1.Activate()
$ws1.Range("A4").Select()
